# Update vm_pu.xlsx res_bus values for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040937264949131
$ws.Cells.Item(2, 4).Value = 1.049243040488897
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.056369171670514
$ws.Cells.Item(2, 9).Value = 1.043919016387306
$ws.Cells.Item(2, 10).Value = 1.046021175515715
$ws.Cells.Item(2, 11).Value = 1.052000642009906
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.059107071904279
$ws.Cells.Item(2, 14).Value = 1.047506645565254
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.041808642010865
$ws.Cells.Item(3, 4).Value = 1.049957520938413
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.057258011100018
$ws.Cells.Item(3, 9).Value = 1.044163549028801
$ws.Cells.Item(3, 10).Value = 1.046538764985156
$ws.Cells.Item(3, 11).Value = 1.052527581964361
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.059809345194445
$ws.Cells.Item(3, 14).Value = 1.048024970071109
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042372748091604
$ws.Cells.Item(4, 4).Value = 1.050420022370961
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.057833799050167
$ws.Cells.Item(4, 9).Value = 1.04432054626854
$ws.Cells.Item(4, 10).Value = 1.046873288506397
$ws.Cells.Item(4, 11).Value = 1.05286805906132
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.060263773983092
$ws.Cells.Item(4, 14).Value = 1.048359968654121
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042609960909866
$ws.Cells.Item(5, 4).Value = 1.050614500703997
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.058076014191845
$ws.Cells.Item(5, 9).Value = 1.044386252336392
$ws.Cells.Item(5, 10).Value = 1.047013827426154
$ws.Cells.Item(5, 11).Value = 1.053011077594502
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.060454817126871
$ws.Cells.Item(5, 14).Value = 1.048500707155263
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042649793621321
$ws.Cells.Item(6, 4).Value = 1.050647156933611
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.058116692157982
$ws.Cells.Item(6, 9).Value = 1.044397267321704
$ws.Cells.Item(6, 10).Value = 1.047037418941416
$ws.Cells.Item(6, 11).Value = 1.05303508408052
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.060486894147062
$ws.Cells.Item(6, 14).Value = 1.048524332173184
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.042375917496408
$ws.Cells.Item(7, 4).Value = 1.050422620833964
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.057837034938057
$ws.Cells.Item(7, 9).Value = 1.044321425398495
$ws.Cells.Item(7, 10).Value = 1.046875166768012
$ws.Cells.Item(7, 11).Value = 1.052869970547202
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.060266326707596
$ws.Cells.Item(7, 14).Value = 1.048361849583082
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041231694892895
$ws.Cells.Item(8, 4).Value = 1.049484463328362
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.05666942370339
$ws.Cells.Item(8, 9).Value = 1.044001911959717
$ws.Cells.Item(8, 10).Value = 1.046196177727747
$ws.Cells.Item(8, 11).Value = 1.052178824252453
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.059344405356711
$ws.Cells.Item(8, 14).Value = 1.047681896300499
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039217528006513
$ws.Cells.Item(9, 4).Value = 1.04783279200302
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.054616986954327
$ws.Cells.Item(9, 9).Value = 1.043429485440211
$ws.Cells.Item(9, 10).Value = 1.044996760278178
$ws.Cells.Item(9, 11).Value = 1.050957246771624
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.05772000622889
$ws.Cells.Item(9, 14).Value = 1.046480775540577
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037876235828152
$ws.Cells.Item(10, 4).Value = 1.046732763569756
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.053252180045797
$ws.Cells.Item(10, 9).Value = 1.04304158649706
$ws.Cells.Item(10, 10).Value = 1.0441952256103
$ws.Cells.Item(10, 11).Value = 1.050140448480091
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.056637250252429
$ws.Cells.Item(10, 14).Value = 1.045678102601532
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037295809208183
$ws.Cells.Item(11, 4).Value = 1.046256715310129
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.052662049466073
$ws.Cells.Item(11, 9).Value = 1.042872140592187
$ws.Cells.Item(11, 10).Value = 1.043847708953515
$ws.Cells.Item(11, 11).Value = 1.049786207387523
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.056168463526442
$ws.Cells.Item(11, 14).Value = 1.045330092431235
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.037080268457163
$ws.Cells.Item(12, 4).Value = 1.046079932035786
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.05244297641147
$ws.Cells.Item(12, 9).Value = 1.042808978666787
$ws.Cells.Item(12, 10).Value = 1.043718559662547
$ws.Cells.Item(12, 11).Value = 1.049654543130035
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.055994344652809
$ws.Cells.Item(12, 14).Value = 1.04520075973346
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.037126500147521
$ws.Cells.Item(13, 4).Value = 1.046117850712748
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.052489962521024
$ws.Cells.Item(13, 9).Value = 1.042822537157526
$ws.Cells.Item(13, 10).Value = 1.043746265600345
$ws.Cells.Item(13, 11).Value = 1.049682789314052
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.056031693296316
$ws.Cells.Item(13, 14).Value = 1.045228505016868
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037277991402586
$ws.Cells.Item(14, 4).Value = 1.046242101479157
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.052643938204691
$ws.Cells.Item(14, 9).Value = 1.042866924136216
$ws.Cells.Item(14, 10).Value = 1.043837034776092
$ws.Cells.Item(14, 11).Value = 1.049775325672756
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.056154070606904
$ws.Cells.Item(14, 14).Value = 1.045319403095256
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037371337559819
$ws.Cells.Item(15, 4).Value = 1.046318662119403
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.052738824675359
$ws.Cells.Item(15, 9).Value = 1.042894242998222
$ws.Cells.Item(15, 10).Value = 1.043892951927992
$ws.Cells.Item(15, 11).Value = 1.049832329366484
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.056229472590332
$ws.Cells.Item(15, 14).Value = 1.045375399655925
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037914764554493
$ws.Cells.Item(16, 4).Value = 1.046764363167848
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.053291362894726
$ws.Cells.Item(16, 9).Value = 1.0430528008717
$ws.Cells.Item(16, 10).Value = 1.044218279817218
$ws.Cells.Item(16, 11).Value = 1.050163946546838
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.056668363328183
$ws.Cells.Item(16, 14).Value = 1.045701189548069
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03825573964063
$ws.Cells.Item(17, 4).Value = 1.047044013346232
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.053638181455989
$ws.Cells.Item(17, 9).Value = 1.04315186342491
$ws.Cells.Item(17, 10).Value = 1.044422230479594
$ws.Cells.Item(17, 11).Value = 1.050371811543088
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.05694368321188
$ws.Cells.Item(17, 14).Value = 1.045905429843779
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038454659462654
$ws.Cells.Item(18, 4).Value = 1.047207154642181
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.053840555781761
$ws.Cells.Item(18, 9).Value = 1.043209501794802
$ws.Cells.Item(18, 10).Value = 1.044541148267954
$ws.Cells.Item(18, 11).Value = 1.050493001376123
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.057104277740764
$ws.Cells.Item(18, 14).Value = 1.046024516509038
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038522491851809
$ws.Cells.Item(19, 4).Value = 1.047262785986286
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.05390957387658
$ws.Cells.Item(19, 9).Value = 1.043229130690472
$ws.Cells.Item(19, 10).Value = 1.044581688805643
$ws.Cells.Item(19, 11).Value = 1.050534314733322
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.05715903714987
$ws.Cells.Item(19, 14).Value = 1.046065114618942
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038219152630293
$ws.Cells.Item(20, 4).Value = 1.04701400682055
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.053600962724818
$ws.Cells.Item(20, 9).Value = 1.043141249751566
$ws.Cells.Item(20, 10).Value = 1.044400352962986
$ws.Cells.Item(20, 11).Value = 1.050349515205022
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.056914143457071
$ws.Cells.Item(20, 14).Value = 1.045883521258588
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037233379461062
$ws.Cells.Item(21, 4).Value = 1.046205511547234
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.052598592665538
$ws.Cells.Item(21, 9).Value = 1.042853859405129
$ws.Cells.Item(21, 10).Value = 1.04381030733369
$ws.Cells.Item(21, 11).Value = 1.049748078313393
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.056118033269384
$ws.Cells.Item(21, 14).Value = 1.045292637696819
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036613906667855
$ws.Cells.Item(22, 4).Value = 1.045697423631269
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.051969102528305
$ws.Cells.Item(22, 9).Value = 1.042671880591956
$ws.Cells.Item(22, 10).Value = 1.043438940108183
$ws.Cells.Item(22, 11).Value = 1.049369449413902
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.055617542168983
$ws.Cells.Item(22, 14).Value = 1.044920743087256
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036942270003429
$ws.Cells.Item(23, 4).Value = 1.045966746916418
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.052302736529324
$ws.Cells.Item(23, 9).Value = 1.042768472599689
$ws.Cells.Item(23, 10).Value = 1.043635844786367
$ws.Cells.Item(23, 11).Value = 1.049570213046444
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.055882856321967
$ws.Cells.Item(23, 14).Value = 1.045117927392668
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038235684617901
$ws.Cells.Item(24, 4).Value = 1.047027565397476
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.053617780017487
$ws.Cells.Item(24, 9).Value = 1.043146046056031
$ws.Cells.Item(24, 10).Value = 1.044410238604559
$ws.Cells.Item(24, 11).Value = 1.050359590128355
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.056927491185126
$ws.Cells.Item(24, 14).Value = 1.045893420938906
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.039737981986605
$ws.Cells.Item(25, 4).Value = 1.048259603785223
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.055146983757024
$ws.Cells.Item(25, 9).Value = 1.043578581521265
$ws.Cells.Item(25, 10).Value = 1.04530718190413
$ws.Cells.Item(25, 11).Value = 1.051273484112622
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.058139927029107
$ws.Cells.Item(25, 14).Value = 1.046791638000843
